# Test script for Invalid login
$wb = $excel.ActiveWorkbook

# Add the new "InvalidLogin" worksheet after the existing ValidLogin sheet.
$validLogin = $wb.Worksheets.Item("ValidLogin")
$newSheet = $wb.Worksheets.Add($null, $validLogin)
$newSheet.Name = "InvalidLogin"

# Populate header row + invalid credentials row.
$newSheet.Range("A1").Value = "Username"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "abcd"
$newSheet.Range("B2").Value = "xyz"

# Selection on the ValidLogin sheet becomes A1:B2 (no active tab there anymore).
$validLogin.Range("A1:B2").Select() | Out-Null

# Make InvalidLogin the active / selected sheet with B3 selected.
$newSheet.Activate() | Out-Null
$newSheet.Range("B3").Select() | Out-Null
